# Applies the updated cryptocurrency price/volume snapshot to Sheet1.
# Cells in columns D (Price) and E (Volume 1h) hold text-formatted values
# (dotted thousand separators, padded percentages); a handful of Price
# values look like plain decimals and would otherwise be auto-coerced to
# numbers by Excel, so those cells are forced to Text format first.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.862.07"
$ws.Range("E2").Value = "  -0.36%  "
$ws.Range("D3").Value = "1.871.04"
$ws.Range("E3").Value = "  -1.26%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "0.7336"
$ws.Range("E5").Value = "  -5.05%  "
$ws.Range("D6").Value = "241.88"
$ws.Range("E6").Value = "  -1.17%  "
$ws.Range("D7").Value = "0.9999"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "0.3152"
$ws.Range("E8").Value = "  +0.46%  "
$ws.Range("D9").Value = "24.61"
$ws.Range("E9").Value = "  -4.63%  "
$ws.Range("D10").Value = "0.07089"
$ws.Range("E10").Value = "  -2.46%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08470"
$ws.Range("E11").Value = "  +2.37%  "
$ws.Range("D12").Value = "0.7498"
$ws.Range("E12").Value = "  -3.12%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.873.68"
$ws.Range("E13").Value = "  -1.59%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "5.367"
$ws.Range("E14").Value = "  -1.95%  "
$ws.Range("D15").Value = "92.39"
$ws.Range("E15").Value = "  -2.77%  "
$ws.Range("D16").Value = "29.864.75"
$ws.Range("E16").Value = "  -0.34%  "
$ws.Range("D17").Value = "6.033"
$ws.Range("E17").Value = "  -2.72%  "
$ws.Range("E18").Value = "  -3.14%  "
$ws.Range("D19").Value = "242.85"
$ws.Range("E19").Value = "  -1.78%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007802"
$ws.Range("E20").Value = "  -0.81%  "
$ws.Range("D21").Value = "0.9989"
$ws.Range("E21").Value = "  -0.11%  "
$ws.Range("D22").Value = "2.122.20"
$ws.Range("E22").Value = "  -0.94%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.910"
$ws.Range("E23").Value = "  -2.94%  "
$ws.Range("D24").Value = "1.001"
$ws.Range("E24").Value = "  +0.06%  "
$ws.Range("D25").Value = "0.1561"
$ws.Range("E25").Value = "  -2.43%  "
$ws.Range("D26").Value = "9.311"
$ws.Range("E26").Value = "  -2.52%  "
$ws.Range("D27").Value = "164.01"
$ws.Range("E27").Value = "  +0.74%  "
$ws.Range("D28").Value = "18.59"
$ws.Range("E28").Value = "  -1.08%  "
$ws.Range("D29").Value = "2.021"
$ws.Range("E29").Value = "  -0.87%  "
$ws.Range("D30").Value = "1.465"
$ws.Range("E30").Value = "  +3.22%  "
$ws.Range("D31").Value = "4.547"
$ws.Range("E31").Value = "  -0.24%  "
$ws.Range("E32").Value = "  -1.43%  "
$ws.Range("D33").Value = "4.261"
$ws.Range("E33").Value = "  +3.68%  "
$ws.Range("D34").Value = "0.05323"
$ws.Range("E34").Value = "  -2.64%  "
$ws.Range("D35").Value = "1.233"
$ws.Range("E35").Value = "  -1.24%  "
$ws.Range("D36").Value = "0.7493"
$ws.Range("E36").Value = "  -0.36%  "
$ws.Range("D37").Value = "0.9991"
$ws.Range("E37").Value = "  -0.39%  "
$ws.Range("D38").Value = "2.697"
$ws.Range("E38").Value = "  +0.87%  "
$ws.Range("D39").Value = "0.01947"
$ws.Range("E39").Value = "  +0.50%  "
$ws.Range("D40").Value = "2.751"
$ws.Range("E40").Value = "  -1.42%  "
$ws.Range("D41").Value = "0.4458"
$ws.Range("E41").Value = "  -0.82%  "
$ws.Range("D42").Value = "1.102.28"
$ws.Range("E42").Value = "  +0.67%  "
$ws.Range("D43").Value = "6.064"
$ws.Range("E43").Value = "  -0.25%  "
$ws.Range("D44").Value = "72.22"
$ws.Range("E44").Value = "  -2.95%  "
$ws.Range("D45").Value = "0.8651"
$ws.Range("E45").Value = "  +1.62%  "
$ws.Range("E46").Value = "  +0.08%  "
$ws.Range("B47").Value = "Quant"
$ws.Range("C47").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "102.70"
$ws.Range("E47").Value = "  -0.05%  "
$ws.Range("B48").Value = "Aptos"
$ws.Range("C48").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D48").Value = "7.712"
$ws.Range("E48").Value = "  +1.32%  "
$ws.Range("D49").Value = "3.072"
$ws.Range("E49").Value = "  +1.98%  "
$ws.Range("D50").Value = "1.835"
$ws.Range("E50").Value = "  -3.16%  "
$ws.Range("D51").Value = "2.020.80"
$ws.Range("E51").Value = "  -0.42%  "
